$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 148 already exists (date "02-08-2021" in A148); this update just fills in
# the MOVE figure (B148) and revises the VIX figure (C148).
$ws.Range("B148").Value = 64.29000000000001
$ws.Range("C148").Value = 19.46

# Daily MOVE/VIX figures for 03-08-2021 .. 03-09-2021, appended as new rows
# 149-172 (Serie date in column A, MOVE in B, VIX in C).
$newRows = @(
    @{ Row = 149; Serie = "03-08-2021"; Move = 65.42; Vix = 18.04 },
    @{ Row = 150; Serie = "04-08-2021"; Move = 62.67; Vix = 17.97 },
    @{ Row = 151; Serie = "05-08-2021"; Move = 63.57; Vix = 17.28 },
    @{ Row = 152; Serie = "06-08-2021"; Move = 62.64; Vix = 16.15 },
    @{ Row = 153; Serie = "09-08-2021"; Move = 65.72; Vix = 16.72 },
    @{ Row = 154; Serie = "10-08-2021"; Move = 66.55; Vix = 16.79 },
    @{ Row = 155; Serie = "11-08-2021"; Move = 60.66; Vix = 16.06 },
    @{ Row = 156; Serie = "12-08-2021"; Move = 57.39; Vix = 15.59 },
    @{ Row = 157; Serie = "13-08-2021"; Move = 55.45; Vix = 15.45 },
    @{ Row = 158; Serie = "16-08-2021"; Move = 59.4; Vix = 16.12 },
    @{ Row = 159; Serie = "17-08-2021"; Move = 60.28; Vix = 17.91 },
    @{ Row = 160; Serie = "18-08-2021"; Move = 58.05; Vix = 21.57 },
    @{ Row = 161; Serie = "19-08-2021"; Move = 59.38; Vix = 21.67 },
    @{ Row = 162; Serie = "20-08-2021"; Move = 59.95; Vix = 18.56 },
    @{ Row = 163; Serie = "23-08-2021"; Move = 62.23; Vix = 17.15 },
    @{ Row = 164; Serie = "24-08-2021"; Move = 63.03; Vix = 17.22 },
    @{ Row = 165; Serie = "25-08-2021"; Move = 64.51000000000001; Vix = 16.79 },
    @{ Row = 166; Serie = "26-08-2021"; Move = 63.41; Vix = 18.84 },
    @{ Row = 167; Serie = "27-08-2021"; Move = 57.98; Vix = 16.39 },
    @{ Row = 168; Serie = "30-08-2021"; Move = 60.14; Vix = 16.19 },
    @{ Row = 169; Serie = "31-08-2021"; Move = 59.54; Vix = 16.48 },
    @{ Row = 170; Serie = "01-09-2021"; Move = 58.09; Vix = 16.11 },
    @{ Row = 171; Serie = "02-09-2021"; Move = 55.8; Vix = 16.41 },
    @{ Row = 172; Serie = "03-09-2021"; Move = 53.26; Vix = 16.41 }
)

foreach ($item in $newRows) {
    # Column A holds text labels like "03-08-2021" elsewhere in the sheet (shared
    # strings, no cell style). Excel normally auto-converts a dd-mm-yyyy-looking
    # string typed into a General cell into a date serial, so the cell is forced
    # to Text first; the temporary formatting is cleared again afterwards so the
    # cell ends up plain/unstyled, matching its neighbours.
    $a = $ws.Range("A" + $item.Row)
    $a.NumberFormat = "@"
    $a.Value = $item.Serie
    $a.ClearFormats()

    $ws.Range("B" + $item.Row).Value = $item.Move
    $ws.Range("C" + $item.Row).Value = $item.Vix
}
